# feat: add 2022-Q1 data
#
#  - rename the existing "总计" sheet to "2022-Q1" and replace its contents
#    with the new per-fund holdings table for that quarter (columns grow
#    from B:D to B:H)
#  - duplicate that sheet (so the new one starts out with the same
#    sheet-level setup: margins, outline settings, etc.) right after
#    itself, rename the duplicate back to "总计", and replace its contents
#    with the updated rollup table (now including the 2022-Q1 row)

$wb = $excel.ActiveWorkbook

# ---- 1. Repurpose the old "总计" sheet as the new "2022-Q1" sheet ----
$q1 = $wb.Worksheets.Item("总计")
$q1.Name = "2022-Q1"

# header row (keeps the existing header style already on B1:D1; just change
# the text and extend it across the two new columns)
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# carry the header formatting (bold/border/center) from the existing header
# cells onto the two newly-used header columns
$q1.Range("D1").Copy() | Out-Null
$q1.Range("E1:H1").PasteSpecial(-4122) | Out-Null

# fund codes / numeric-looking figures are stored as text in this data set,
# so format those columns as Text before writing so Excel doesn't coerce
# them to numbers (and drop the leading zero on the fund codes)
$q1.Range("B2:B3").NumberFormat = "@"
$q1.Range("D2:G3").NumberFormat = "@"

# row 2
$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "011179"
$q1.Range("C2").Value = "浙商智选食品饮料股票A"
$q1.Range("D2").Value = "0.22"
$q1.Range("E2").Value = "91.35"
$q1.Range("F2").Value = "5.65"
$q1.Range("G2").Value = "0.0124"
$q1.Range("H2").Value = 9

# row 3
$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "011180"
$q1.Range("C3").Value = "浙商智选食品饮料股票C"
$q1.Range("D3").Value = "0.05"
$q1.Range("E3").Value = "91.35"
$q1.Range("F3").Value = "5.65"
$q1.Range("G3").Value = "0.0028"
$q1.Range("H3").Value = 9

# ---- 2. Duplicate "2022-Q1" right after itself to become the new "总计" ----
$q1.Copy($null, $q1)
$total = $wb.Worksheets.Item($q1.Index + 1)
$total.Name = "总计"

# wipe the duplicated fund table - we only want the rollup columns back
$total.Cells.Clear()

# re-apply the header style (s=2) and the row-index style from "2022-Q1"
# onto the now-empty cells before writing the rollup values
$q1.Range("B1").Copy() | Out-Null
$total.Range("B1:D1").PasteSpecial(-4122) | Out-Null
$q1.Range("A2").Copy() | Out-Null
$total.Range("A2:A4").PasteSpecial(-4122) | Out-Null

$total.Range("B1").Value = "日期"
$total.Range("C1").Value = "持有数量(只)"
$total.Range("D1").Value = "持有市值(亿元)"

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.02

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2021-Q4"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.01

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2021-Q3"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.02
